$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap full row content (columns B:AB) between paired rows ---
# These pairs represent matches whose data rows were reordered in the source export
$swapPairs = @(
    @(2, 3),
    @(7, 8),
    @(19, 20),
    @(31, 32),
    @(38, 39),
    @(43, 44),
    @(50, 51),
    @(124, 125),
    @(140, 141),
    @(190, 191)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# --- Append new match rows (266-273) at the end of the table ---
# Row 266
$ws.Range("A265").Copy($ws.Range("A266"))
$ws.Range("D265").Copy($ws.Range("D266"))
$ws.Range("A266").Value = 264
$ws.Range("B266").Value = 8106962
$ws.Range("C266").Value = "Serbia Super Liga"
$ws.Range("D266").Value = 45419.375
$ws.Range("E266").Value = "IMT Novi Belgrade"
$ws.Range("F266").Value = "FK Novi Pazar"
$ws.Range("G266").Value = 0
$ws.Range("H266").Value = 1
$ws.Range("I266").Value = "A"
$ws.Range("J266").Value = 2.4
$ws.Range("K266").Value = 3.3
$ws.Range("L266").Value = 2.5
$ws.Range("M266").Value = 1.571
$ws.Range("N266").Value = 4.2
$ws.Range("O266").Value = 4
$ws.Range("P266").Value = -1
$ws.Range("Q266").Value = 1.975
$ws.Range("R266").Value = 1.825
$ws.Range("S266").Value = 3
$ws.Range("T266").Value = 2.025
$ws.Range("U266").Value = 1.775
$ws.Range("V266").Value = -1
$ws.Range("W266").Value = -1
$ws.Range("X266").Value = 3
$ws.Range("Y266").Value = -1
$ws.Range("Z266").Value = 0.825
$ws.Range("AA266").Value = -1
$ws.Range("AB266").Value = 0.7749999999999999

# Row 267
$ws.Range("A265").Copy($ws.Range("A267"))
$ws.Range("D265").Copy($ws.Range("D267"))
$ws.Range("A267").Value = 265
$ws.Range("B267").Value = 8105038
$ws.Range("C267").Value = "Serbia Super Liga"
$ws.Range("D267").Value = 45419.54166666666
$ws.Range("E267").Value = "FK Radnik Surdulica"
$ws.Range("F267").Value = "Javor Ivanjica"
$ws.Range("G267").Value = 0
$ws.Range("H267").Value = 1
$ws.Range("I267").Value = "A"
$ws.Range("J267").Value = 2.5
$ws.Range("K267").Value = 3.2
$ws.Range("L267").Value = 2.5
$ws.Range("M267").Value = 3.5
$ws.Range("N267").Value = 3.3
$ws.Range("O267").Value = 1.85
$ws.Range("P267").Value = 0.5
$ws.Range("Q267").Value = 1.875
$ws.Range("R267").Value = 1.925
$ws.Range("S267").Value = 2.25
$ws.Range("T267").Value = 1.775
$ws.Range("U267").Value = 2.025
$ws.Range("V267").Value = -1
$ws.Range("W267").Value = -1
$ws.Range("X267").Value = 0.8500000000000001
$ws.Range("Y267").Value = -1
$ws.Range("Z267").Value = 0.925
$ws.Range("AA267").Value = -1
$ws.Range("AB267").Value = 1.025

# Row 268
$ws.Range("A265").Copy($ws.Range("A268"))
$ws.Range("D265").Copy($ws.Range("D268"))
$ws.Range("A268").Value = 266
$ws.Range("B268").Value = 8106961
$ws.Range("C268").Value = "Serbia Super Liga"
$ws.Range("D268").Value = 45419.54166666666
$ws.Range("E268").Value = "Spartak Subotica"
$ws.Range("F268").Value = "Radnicki Nis"
$ws.Range("G268").Value = 3
$ws.Range("H268").Value = 1
$ws.Range("I268").Value = "H"
$ws.Range("J268").Value = 2.4
$ws.Range("K268").Value = 3.2
$ws.Range("L268").Value = 2.625
$ws.Range("M268").Value = 3
$ws.Range("N268").Value = 2.15
$ws.Range("O268").Value = 3.3
$ws.Range("P268").Value = 0
$ws.Range("Q268").Value = 1.8
$ws.Range("R268").Value = 2
$ws.Range("S268").Value = 2
$ws.Range("T268").Value = 1.85
$ws.Range("U268").Value = 1.95
$ws.Range("V268").Value = 2
$ws.Range("W268").Value = -1
$ws.Range("X268").Value = -1
$ws.Range("Y268").Value = 0.8
$ws.Range("Z268").Value = -1
$ws.Range("AA268").Value = 0.8500000000000001
$ws.Range("AB268").Value = -1

# Row 269
$ws.Range("A265").Copy($ws.Range("A269"))
$ws.Range("D265").Copy($ws.Range("D269"))
$ws.Range("A269").Value = 267
$ws.Range("B269").Value = 8105037
$ws.Range("C269").Value = "Serbia Super Liga"
$ws.Range("D269").Value = 45419.58333333334
$ws.Range("E269").Value = "FK Vozdovac"
$ws.Range("F269").Value = "FK Zeleznicar Pancevo"
$ws.Range("G269").Value = 2
$ws.Range("H269").Value = 3
$ws.Range("I269").Value = "A"
$ws.Range("J269").Value = 2
$ws.Range("K269").Value = 3.25
$ws.Range("L269").Value = 3.3
$ws.Range("M269").Value = 2.45
$ws.Range("N269").Value = 3.25
$ws.Range("O269").Value = 2.55
$ws.Range("P269").Value = 0
$ws.Range("Q269").Value = 1.85
$ws.Range("R269").Value = 1.95
$ws.Range("S269").Value = 2.5
$ws.Range("T269").Value = 1.825
$ws.Range("U269").Value = 1.975
$ws.Range("V269").Value = -1
$ws.Range("W269").Value = -1
$ws.Range("X269").Value = 1.55
$ws.Range("Y269").Value = -1
$ws.Range("Z269").Value = 0.95
$ws.Range("AA269").Value = 0.825
$ws.Range("AB269").Value = -1

# Row 270
$ws.Range("A265").Copy($ws.Range("A270"))
$ws.Range("D265").Copy($ws.Range("D270"))
$ws.Range("A270").Value = 268
$ws.Range("B270").Value = 8105864
$ws.Range("C270").Value = "Serbia Super Liga"
$ws.Range("D270").Value = 45420.45833333334
$ws.Range("E270").Value = "FK Backa Topola"
$ws.Range("F270").Value = "FK Radnicki 1923"
$ws.Range("G270").Value = 4
$ws.Range("H270").Value = 3
$ws.Range("I270").Value = "H"
$ws.Range("J270").Value = 1.533
$ws.Range("K270").Value = 4
$ws.Range("L270").Value = 4.75
$ws.Range("M270").Value = 1.615
$ws.Range("N270").Value = 4
$ws.Range("O270").Value = 4.2
$ws.Range("P270").Value = -0.75
$ws.Range("Q270").Value = 1.8
$ws.Range("R270").Value = 2
$ws.Range("S270").Value = 3
$ws.Range("T270").Value = 1.825
$ws.Range("U270").Value = 1.975
$ws.Range("V270").Value = 0.615
$ws.Range("W270").Value = -1
$ws.Range("X270").Value = -1
$ws.Range("Y270").Value = 0.4
$ws.Range("Z270").Value = -0.5
$ws.Range("AA270").Value = 0.825
$ws.Range("AB270").Value = -1

# Row 271
$ws.Range("A265").Copy($ws.Range("A271"))
$ws.Range("D265").Copy($ws.Range("D271"))
$ws.Range("A271").Value = 269
$ws.Range("B271").Value = 8105018
$ws.Range("C271").Value = "Serbia Super Liga"
$ws.Range("D271").Value = 45420.45833333334
$ws.Range("E271").Value = "FK Napredak"
$ws.Range("F271").Value = "Partizan Belgrade"
$ws.Range("G271").Value = 3
$ws.Range("H271").Value = 3
$ws.Range("I271").Value = "D"
$ws.Range("J271").Value = 6
$ws.Range("K271").Value = 4.6
$ws.Range("L271").Value = 1.363
$ws.Range("M271").Value = 9
$ws.Range("N271").Value = 5.75
$ws.Range("O271").Value = 1.25
$ws.Range("P271").Value = 1.75
$ws.Range("Q271").Value = 1.925
$ws.Range("R271").Value = 1.875
$ws.Range("S271").Value = 3.25
$ws.Range("T271").Value = 1.975
$ws.Range("U271").Value = 1.825
$ws.Range("V271").Value = -1
$ws.Range("W271").Value = 4.75
$ws.Range("X271").Value = -1
$ws.Range("Y271").Value = 0.925
$ws.Range("Z271").Value = -1
$ws.Range("AA271").Value = 0.9750000000000001
$ws.Range("AB271").Value = -1

# Row 272
$ws.Range("A265").Copy($ws.Range("A272"))
$ws.Range("D265").Copy($ws.Range("D272"))
$ws.Range("A272").Value = 270
$ws.Range("B272").Value = 8105019
$ws.Range("C272").Value = "Serbia Super Liga"
$ws.Range("D272").Value = 45420.54166666666
$ws.Range("E272").Value = "Vojvodina"
$ws.Range("F272").Value = "FK Cukaricki"
$ws.Range("G272").Value = 2
$ws.Range("H272").Value = 3
$ws.Range("I272").Value = "A"
$ws.Range("J272").Value = 2.1
$ws.Range("K272").Value = 3.4
$ws.Range("L272").Value = 2.9
$ws.Range("M272").Value = 2
$ws.Range("N272").Value = 3.6
$ws.Range("O272").Value = 2.9
$ws.Range("P272").Value = -0.25
$ws.Range("Q272").Value = 1.825
$ws.Range("R272").Value = 1.975
$ws.Range("S272").Value = 2.75
$ws.Range("T272").Value = 1.8
$ws.Range("U272").Value = 2
$ws.Range("V272").Value = -1
$ws.Range("W272").Value = -1
$ws.Range("X272").Value = 1.9
$ws.Range("Y272").Value = -1
$ws.Range("Z272").Value = 0.9750000000000001
$ws.Range("AA272").Value = 0.8
$ws.Range("AB272").Value = -1

# Row 273
$ws.Range("A265").Copy($ws.Range("A273"))
$ws.Range("D265").Copy($ws.Range("D273"))
$ws.Range("A273").Value = 271
$ws.Range("B273").Value = 8105020
$ws.Range("C273").Value = "Serbia Super Liga"
$ws.Range("D273").Value = 45420.5625
$ws.Range("E273").Value = "Mladost Lucani"
$ws.Range("F273").Value = "Crvena Zvezda"
$ws.Range("G273").Value = 0
$ws.Range("H273").Value = 1
$ws.Range("I273").Value = "A"
$ws.Range("J273").Value = 10
$ws.Range("K273").Value = 6.5
$ws.Range("L273").Value = 1.166
$ws.Range("M273").Value = 9
$ws.Range("N273").Value = 6.5
$ws.Range("O273").Value = 1.2
$ws.Range("P273").Value = 2
$ws.Range("Q273").Value = 1.9
$ws.Range("R273").Value = 1.9
$ws.Range("S273").Value = 3.5
$ws.Range("T273").Value = 1.9
$ws.Range("U273").Value = 1.9
$ws.Range("V273").Value = -1
$ws.Range("W273").Value = -1
$ws.Range("X273").Value = 0.2
$ws.Range("Y273").Value = 0.8999999999999999
$ws.Range("Z273").Value = -1
$ws.Range("AA273").Value = -1
$ws.Range("AB273").Value = 0.8999999999999999
